# Add 8 new "Health Care Assistants" COVID resource rows to the
# "COVID Resources-HCP" sheet (rows 114-121), matching the pattern of the
# existing rows (A=Audience, B=Topic, C=Location, D=Name, E=Author,
# F=Resource Type, G=Link [hyperlinked]).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID Resources-HCP")

$rows = @(
    @{ Row=114; A="Health Care Assistants"; B="Healthcare Provider Wellness"; C="British Columbia"; D="Health Care Assistants Summer 2020 Self-Reflection to Support Grief and Growth"; E="BC Centre for Palliative Care"; F="Education"; G="https://www.youtube.com/watch?v=eawdB0kqx60&feature=youtu.be" },
    @{ Row=115; A="Health Care Assistants"; B="Psychosocial Care";            C="United States";    D="The Gift and Power of Emotional Courage";                                                    E="Susan David";                  F="Education"; G="https://www.ted.com/talks/susan_david_the_gift_and_power_of_emotional_courage" },
    @{ Row=116; A="Health Care Assistants"; B="Psychosocial Care";            C="United States";    D="How Journaling Can Help You in Hard Times";                                                  E="Kira M. Newman";               F="Education"; G="https://greatergood.berkeley.edu/article/item/how_journaling_can_help_you_in_hard_times" },
    @{ Row=117; A="Health Care Assistants"; B="Healthcare Provider Wellness"; C="British Columbia"; D="Health Care Assistant Summer 2020 Adapting Rituals of Grief & Growth";                       E="BC Centre for Palliative Care"; F="Education"; G="https://www.youtube.com/watch?v=XXb4j8gRxEw&t=6s" },
    @{ Row=118; A="Health Care Assistants"; B="Healthcare Provider Wellness"; C="British Columbia"; D="Health Care Assistant Summer 2020 Social Connection in a Time of Physical Distancing";        E="BC Centre for Palliative Care"; F="Education"; G="https://bc-cpc.ca/cpc/wp-content/uploads/2020/08/BCCPC-Patio-Ponderings-1-Social-Connection.pdf" },
    @{ Row=119; A="Health Care Assistants"; B="Healthcare Provider Wellness"; C="British Columbia"; D="Health Care Assistant Summer 2020 Dealing with Fear & Anxiety";                               E="BC Centre for Palliative Care"; F="Education"; G="https://bc-cpc.ca/cpc/wp-content/uploads/2020/08/BCCPC-Patio-Ponderings-2-Fear-and-Anxiety.pdf" },
    @{ Row=120; A="Health Care Assistants"; B="Healthcare Provider Wellness"; C="British Columbia"; D="Health Care Assistant Summer 2020 Adapting Rituals of Grief & Growth";                       E="BC Centre for Palliative Care"; F="Education"; G="https://bc-cpc.ca/cpc/wp-content/uploads/2020/08/BCCPC-Patio-Ponderings-3-Rituals.pdf" },
    @{ Row=121; A="Health Care Assistants"; B="Healthcare Provider Wellness"; C="British Columbia"; D="Health Care Assistant Summer 2020 Self-Reflection to Support Grief and Growth";                E="BC Centre for Palliative Care"; F="Education"; G="https://bc-cpc.ca/cpc/wp-content/uploads/2020/08/BCCPC-Patio-Ponderings-4-Reflection.pdf" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F
    $ws.Range("G$rowNum").Value = $r.G

    # Match the look of the rows above: left-aligned + wrapped text, row
    # height 30, and a real hyperlink on the Link column.
    $ws.Range("A$rowNum").HorizontalAlignment = -4131  # xlLeft
    $ws.Range("A$rowNum").WrapText = $true
    $ws.Range("C$rowNum").HorizontalAlignment = -4131
    $ws.Range("E$rowNum").HorizontalAlignment = -4131
    $ws.Range("F$rowNum").HorizontalAlignment = -4131
    $ws.Rows.Item($rowNum).RowHeight = 30

    $ws.Hyperlinks.Add($ws.Range("G$rowNum"), $r.G) | Out-Null
    $ws.Range("G$rowNum").Style = "Hyperlink"
}

# Put the view/selection where the author left off after adding the rows.
$ws.Range("G121").Select()

Write-Output "Added rows 114-121 to 'COVID Resources-HCP'"
